# Weekly data refresh: a new week's price record for
# "Macroferia Regional de Talca - Arveja Verde" is inserted as a new
# row right before the existing row 37, pushing the following rows
# (old 37..99) down by one (new 38..100).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 37 (shifts rows 37:99 -> 38:100).
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row 37 with the new week's record.
$ws.Cells.Item(37, 1).Value  = 5
$ws.Cells.Item(37, 2).Value  = 'Macroferia Regional de Talca'
$ws.Cells.Item(37, 3).Value  = 'Maule'
$ws.Cells.Item(37, 4).Value  = 44645
$ws.Cells.Item(37, 5).Value  = 7
$ws.Cells.Item(37, 6).Value  = 100112022
$ws.Cells.Item(37, 7).Value  = 'Arveja Verde'
$ws.Cells.Item(37, 8).Value  = 'Sin especificar'
$ws.Cells.Item(37, 9).Value  = 'Primera'
$ws.Cells.Item(37, 10).Value = 200
$ws.Cells.Item(37, 11).Value = 25000
$ws.Cells.Item(37, 12).Value = 25000
$ws.Cells.Item(37, 13).Value = 25000
$ws.Cells.Item(37, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(37, 15).Value = 'Carahue'
$ws.Cells.Item(37, 16).Value = 1000
$ws.Cells.Item(37, 17).Value = 25
$ws.Cells.Item(37, 18).Value = 'Hortaliza'

# Give the new date cell the same date-time number format as the rest
# of column D (style index 2 uses numFmtId 165 "YYYY-MM-DD HH:MM:SS").
$ws.Cells.Item(37, 4).NumberFormat = $ws.Cells.Item(38, 4).NumberFormat
